# Updates cryptos list prices/volumes per commit 'Updated cryptos list on Tue Sep 19 06:52:21 UTC 2023 with GitHub Actions'.
# Also swaps the Cosmos/Stellar (rows 27-28) and ImmutableX/ARBITRUM (rows 38-39) entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be stored as text even if it looks numeric,
# then clear the temporary number-format so no residual style is left behind.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# --- Row swaps: Cosmos/Stellar (rows 27<->28), ImmutableX/ARBITRUM (rows 38<->39) ---
$ws.Cells.Item(27, 2).Value = 'Stellar'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 27 4 '0.119'
$ws.Cells.Item(27, 5).Value = '  +0.38%  '

$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 28 4 '7.25'
$ws.Cells.Item(28, 5).Value = '  +0.11%  '

$ws.Cells.Item(38, 2).Value = 'ARBITRUM'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 38 4 '0.835'
$ws.Cells.Item(38, 5).Value = '  +2.06%  '

$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 39 4 '0.532'
$ws.Cells.Item(39, 5).Value = '  -0.07%  '

# --- Remaining Price (D) / Volume(1h) (E) updates ---
$ws.Cells.Item(2, 4).Value = '26.857.94'
$ws.Cells.Item(2, 5).Value = '  -0.02%  '
$ws.Cells.Item(3, 4).Value = '1.638.89'
$ws.Cells.Item(3, 5).Value = '  -0.25%  '
$ws.Cells.Item(4, 5).Value = '  -0.61%  '
Set-TextValue 5 4 '216.88'
$ws.Cells.Item(5, 5).Value = '  -0.81%  '
Set-TextValue 6 4 '0.507'
$ws.Cells.Item(6, 5).Value = '  +2.06%  '
$ws.Cells.Item(7, 5).Value = '  -0.58%  '
$ws.Cells.Item(8, 5).Value = '  +1.24%  '
$ws.Cells.Item(9, 5).Value = '  +0.16%  '
$ws.Cells.Item(10, 5).Value = '  +3.14%  '
Set-TextValue 11 4 '0.0845'
$ws.Cells.Item(11, 5).Value = '  +0.08%  '
$ws.Cells.Item(12, 4).Value = '1.867.63'
$ws.Cells.Item(12, 5).Value = '  -0.25%  '
$ws.Cells.Item(13, 4).Value = '1.629.50'
$ws.Cells.Item(13, 5).Value = '  -0.88%  '
$ws.Cells.Item(14, 5).Value = '  -1.00%  '
$ws.Cells.Item(15, 5).Value = '  +0.50%  '
Set-TextValue 16 4 '67.18'
$ws.Cells.Item(16, 5).Value = '  +2.72%  '
$ws.Cells.Item(17, 4).Value = '26.854.76'
$ws.Cells.Item(17, 5).Value = '  -0.06%  '
$ws.Cells.Item(18, 4).Value = '0.0₃0729'
$ws.Cells.Item(18, 5).Value = '  -0.15%  '
Set-TextValue 19 4 '218.31'
$ws.Cells.Item(19, 5).Value = '  +1.27%  '
$ws.Cells.Item(20, 5).Value = '  -0.63%  '
Set-TextValue 21 4 '6.73'
$ws.Cells.Item(21, 5).Value = '  +1.06%  '
$ws.Cells.Item(22, 5).Value = '  +0.64%  '
$ws.Cells.Item(23, 5).Value = '  +2.74%  '
$ws.Cells.Item(24, 5).Value = '  -0.52%  '
$ws.Cells.Item(25, 5).Value = '  -0.42%  '
$ws.Cells.Item(26, 5).Value = '  -0.50%  '
Set-TextValue 29 4 '15.80'
$ws.Cells.Item(29, 5).Value = '  +0.55%  '
$ws.Cells.Item(30, 5).Value = '  -1.15%  '
Set-TextValue 31 4 '1.19'
$ws.Cells.Item(31, 5).Value = '  -1.03%  '
Set-TextValue 32 4 '3.32'
$ws.Cells.Item(32, 5).Value = '  -1.33%  '
$ws.Cells.Item(33, 5).Value = '  -0.24%  '
Set-TextValue 34 4 '1.56'
$ws.Cells.Item(34, 5).Value = '  +1.18%  '
$ws.Cells.Item(35, 4).Value = '1.265.35'
$ws.Cells.Item(35, 5).Value = '  -1.38%  '
Set-TextValue 36 4 '2.43'
$ws.Cells.Item(36, 5).Value = '  -0.16%  '
$ws.Cells.Item(37, 5).Value = '  +1.91%  '
$ws.Cells.Item(40, 5).Value = '  -0.58%  '
$ws.Cells.Item(41, 5).Value = '  +0.21%  '
Set-TextValue 42 4 '5.37'
$ws.Cells.Item(42, 5).Value = '  +0.42%  '
$ws.Cells.Item(43, 4).Value = '1.778.73'
$ws.Cells.Item(43, 5).Value = '  -0.25%  '
Set-TextValue 44 4 '61.87'
$ws.Cells.Item(44, 5).Value = '  +1.41%  '
$ws.Cells.Item(45, 5).Value = '  -0.26%  '
Set-TextValue 46 4 '91.78'
$ws.Cells.Item(46, 5).Value = '  -1.03%  '
$ws.Cells.Item(47, 5).Value = '  -0.61%  '
$ws.Cells.Item(48, 5).Value = '  -0.08%  '
Set-TextValue 49 4 '0.0512'
$ws.Cells.Item(49, 5).Value = '  -0.74%  '
$ws.Cells.Item(50, 5).Value = '  +0.83%  '
$ws.Cells.Item(51, 5).Value = '  -0.70%  '
